$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 369.75
$ws.Range("I42").Value = 76.333336
$ws.Range("J42").Value = 1250
$ws.Range("K42").Value = 229.000008
$ws.Range("L42").Value = 3750
$ws.Range("M42").Value = 0.9999919999999918
$ws.Range("N42").Value = -4210
$ws.Range("H43").Value = 438.25
$ws.Range("I43").Value = 450.5
$ws.Range("J43").Value = 426
$ws.Range("K43").Value = 450.5
$ws.Range("L43").Value = 426
$ws.Range("M43").Value = -381.5
$ws.Range("N43").Value = -564
$ws.Range("H137").Value = 1805.6451
$ws.Range("I137").Value = 1208.15
$ws.Range("J137").Value = 2892
$ws.Range("K137").Value = 3624.45
$ws.Range("L137").Value = 8676
$ws.Range("M137").Value = -1074.45
$ws.Range("N137").Value = -13776
$ws.Range("H138").Value = 144854.25
$ws.Range("I138").Value = 286342.84
$ws.Range("J138").Value = 3365.6858
$ws.Range("K138").Value = 859028.52
$ws.Range("L138").Value = 10097.0574
$ws.Range("M138").Value = -853888.52
$ws.Range("N138").Value = -20377.0574

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12485.969
$ws.Range("I32").Value = 4025.6904
$ws.Range("J32").Value = 27935.174
$ws.Range("K32").Value = 4025.6904
$ws.Range("L32").Value = 27935.174
$ws.Range("M32").Value = -3738.6904
$ws.Range("N32").Value = -28509.174
$ws.Range("H45").Value = 2207.7727
$ws.Range("I45").Value = 2267.4119
$ws.Range("K45").Value = 2267.4119
$ws.Range("M45").Value = -1890.4119
$ws.Range("H74").Value = 5104109
$ws.Range("I74").Value = 8930631
$ws.Range("J74").Value = 2079.1904
$ws.Range("K74").Value = 8930631
$ws.Range("L74").Value = 2079.1904
$ws.Range("M74").Value = -8929757
$ws.Range("N74").Value = -3827.1904
$ws.Range("H77").Value = 5104109
$ws.Range("I77").Value = 8930631
$ws.Range("J77").Value = 2079.1904
$ws.Range("K77").Value = 44653155
$ws.Range("L77").Value = 10395.952
$ws.Range("M77").Value = -44648787
$ws.Range("N77").Value = -19131.952
$ws.Range("H107").Value = 36500
$ws.Range("J107").Value = 36500
$ws.Range("L107").Value = 36500
$ws.Range("N107").Value = -44180
$ws.Range("H109").Value = 61500
$ws.Range("J109").Value = 61500
$ws.Range("L109").Value = 61500
$ws.Range("N109").Value = -64274
$ws.Range("H112").Value = 12600
$ws.Range("J112").Value = 12600
$ws.Range("L112").Value = 12600
$ws.Range("N112").Value = -15554

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 42000
$ws.Range("J108").Value = 42000
$ws.Range("L108").Value = 42000
$ws.Range("N108").Value = -49680
$ws.Range("H110").Value = 27411.111
$ws.Range("J110").Value = 27411.111
$ws.Range("L110").Value = 27411.111
$ws.Range("N110").Value = -35591.111
$ws.Range("H111").Value = 39750
$ws.Range("J111").Value = 39750
$ws.Range("L111").Value = 39750
$ws.Range("N111").Value = -47930
$ws.Range("H112").Value = 35000
$ws.Range("J112").Value = 35000
$ws.Range("L112").Value = 35000
$ws.Range("N112").Value = -37954
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("H118").Value = 32000
$ws.Range("J118").Value = 32000
$ws.Range("L118").Value = 32000
$ws.Range("N118").Value = -35314
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("H120").Value = 35000
$ws.Range("J120").Value = 35000
$ws.Range("L120").Value = 35000
$ws.Range("N120").Value = -44676
$ws.Range("N116").ClearContents()
$ws.Range("N117").ClearContents()
$ws.Range("N119").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 28000
$ws.Range("J108").Value = 28000
$ws.Range("L108").Value = 28000
$ws.Range("N108").Value = -35680
$ws.Range("H109").Value = 22000
$ws.Range("J109").Value = 22000
$ws.Range("L109").Value = 22000
$ws.Range("N109").Value = -24080
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("H115").Value = 40000
$ws.Range("J115").Value = 40000
$ws.Range("L115").Value = 40000
$ws.Range("N115").Value = -42350
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("H118").Value = 40000
$ws.Range("J118").Value = 40000
$ws.Range("L118").Value = 40000
$ws.Range("N118").Value = -43314
$ws.Range("H119").Value = 41187
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 41187
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 41187
$ws.Range("N119").Value = -50863
$ws.Range("H120").Value = 35000
$ws.Range("J120").Value = 35000
$ws.Range("L120").Value = 35000
$ws.Range("N120").Value = -42258
$ws.Range("H134").Value = 5818.769
$ws.Range("I134").Value = 6461.6
$ws.Range("J134").Value = 3676
$ws.Range("K134").Value = 19384.8
$ws.Range("L134").Value = 11028
$ws.Range("M134").Value = -16849.8
$ws.Range("N134").Value = -16098
$ws.Range("N114").ClearContents()
$ws.Range("N116").ClearContents()
$ws.Range("M119").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6173771.5
$ws.Range("J131").Value = 6803680.5
$ws.Range("L131").Value = 20411041.5
$ws.Range("N131").Value = -20421121.5
$ws.Range("H132").Value = 1325.75
$ws.Range("I132").Value = 701.5
$ws.Range("K132").Value = 6313.5
$ws.Range("M132").Value = -3783.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").Value = 56100
$ws.Range("J125").Value = 56100
$ws.Range("L125").Value = 56100
$ws.Range("N125").Value = -61020

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 389
$ws.Range("I55").Value = 342.4
$ws.Range("J55").Value = 466.66666
$ws.Range("K55").Value = 342.4
$ws.Range("L55").Value = 466.66666
$ws.Range("M55").Value = -169.4
$ws.Range("N55").Value = -812.66666

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 111114264
$ws.Range("I81").Value = 333337100
$ws.Range("J81").Value = 2850
$ws.Range("K81").Value = 666674200
$ws.Range("L81").Value = 5700
$ws.Range("M81").Value = -666673139
$ws.Range("N81").Value = -7822
$ws.Range("H84").Value = 111114264
$ws.Range("I84").Value = 333337100
$ws.Range("J84").Value = 2850
$ws.Range("K84").Value = 3333371000
$ws.Range("L84").Value = 28500
$ws.Range("M84").Value = -3333365696
$ws.Range("N84").Value = -39108
